# Change review rating in book reviews to star rating
# Remove the white (bg1) solid fill from the outlined "5-Point Star 7"
# shape on each slide, turning it into a no-fill (outline only) star.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.Name -eq "5-Point Star 7") {
            $shape.Fill.Visible = 0
        }
    }
}
